$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.235.06'
$ws.Range("E2").Value = '  +1.10%  '

# Row 3
$ws.Range("D3").Value = '1.798.88'
$ws.Range("E3").Value = '  +2.24%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.38'
$ws.Range("E5").Value = '  +0.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9985'
$ws.Range("E6").Value = '  -0.23%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4776'
$ws.Range("E7").Value = '  +27.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3640'
$ws.Range("E8").Value = '  +8.95%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.44'
$ws.Range("E9").Value = '  -0.60%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07663'
$ws.Range("E10").Value = '  +7.38%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.143'
$ws.Range("E11").Value = '  +2.28%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.54'
$ws.Range("E12").Value = '  +1.39%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9998'
$ws.Range("E13").Value = '  -0.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.286'
$ws.Range("E14").Value = '  +1.94%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.289'
$ws.Range("E15").Value = '  +2.15%  '

# Row 16
$ws.Range("D16").Value = '1.796.92'
$ws.Range("E16").Value = '  +2.26%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001093'
$ws.Range("E17").Value = '  +4.25%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06711'
$ws.Range("E18").Value = '  +2.19%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.66'
$ws.Range("E19").Value = '  +2.11%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9990'

# Row 21
$ws.Range("E21").Value = '  +2.62%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.402'
$ws.Range("E22").Value = '  +2.49%  '

# Row 23
$ws.Range("D23").Value = '28.252.13'
$ws.Range("E23").Value = '  +1.16%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.00'
$ws.Range("E24").Value = '  +3.03%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.397'
$ws.Range("E25").Value = '  +0.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.56'
$ws.Range("E26").Value = '  +4.40%  '

# Row 27
$ws.Range("E27").Value = '  +4.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.25'
$ws.Range("E28").Value = '  -0.52%  '

# Row 29
$ws.Range("D29").Value = '2.003.69'
$ws.Range("E29").Value = '  +2.35%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.48'
$ws.Range("E30").Value = '  +1.82%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.272'
$ws.Range("E31").Value = '  +0.36%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.038'
$ws.Range("E32").Value = '  +0.43%  '

# Row 33
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09602'
$ws.Range("E33").Value = '  +10.02%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.916'
$ws.Range("E34").Value = '  +2.77%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02387'
$ws.Range("E35").Value = '  +2.58%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.15'
$ws.Range("E36").Value = '  +0.19%  '

# Row 37
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6659'
$ws.Range("E37").Value = '  +1.93%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06292'
$ws.Range("E38").Value = '  +2.06%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.202'
$ws.Range("E39").Value = '  +1.54%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2169'
$ws.Range("E40").Value = '  +3.29%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.480'
$ws.Range("E41").Value = '  +2.04%  '

# Row 42
$ws.Range("E42").Value = '  +0.57%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.073'
$ws.Range("E43").Value = '  +0.88%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.14'
$ws.Range("E44").Value = '  +3.53%  '

# Row 45
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9985'
$ws.Range("E45").Value = '  -0.20%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.869'
$ws.Range("E46").Value = '  +1.08%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6121'
$ws.Range("E47").Value = '  +2.11%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.23'
$ws.Range("E48").Value = '  -0.47%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.037'
$ws.Range("E49").Value = '  +1.92%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.171'
$ws.Range("E50").Value = '  -0.28%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07098'
$ws.Range("E51").Value = '  -0.69%  '
